# Auto-generated edit script: updates market-price derived columns (H-N)
# in the per-job leve profit sheets, per the scheduled price-refresh diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 1849
$ws.Range("I12").Value = 2633.5
$ws.Range("J12").Value = 280
$ws.Range("K12").Value = 2633.5
$ws.Range("L12").Value = 280
$ws.Range("M12").Value = -2463.5
$ws.Range("H48").Value = 1500
$ws.Range("I48").Value = 0
$ws.Range("J48").Value = 1500
$ws.Range("K48").Value = 0
$ws.Range("L48").Value = 4500
$ws.Range("N48").Value = -5084
$ws.Range("M48").ClearContents()
$ws.Range("H56").Value = 1500
$ws.Range("I56").Value = 0
$ws.Range("J56").Value = 1500
$ws.Range("K56").Value = 0
$ws.Range("L56").Value = 4500
$ws.Range("N56").Value = -5568
$ws.Range("M56").ClearContents()
$ws.Range("H86").Value = 7091.3
$ws.Range("I86").Value = 6160
$ws.Range("J86").Value = 8022.6
$ws.Range("K86").Value = 6160
$ws.Range("L86").Value = 8022.6
$ws.Range("M86").Value = -5037
$ws.Range("N86").Value = -10268.6
$ws.Range("H89").Value = 7091.3
$ws.Range("I89").Value = 6160
$ws.Range("J89").Value = 8022.6
$ws.Range("K89").Value = 30800
$ws.Range("L89").Value = 40113
$ws.Range("M89").Value = -25184
$ws.Range("N89").Value = -51345
$ws.Range("H113").Value = 8206.5
$ws.Range("I113").Value = 7819.222
$ws.Range("J113").Value = 9368.333000000001
$ws.Range("K113").Value = 7819.222
$ws.Range("L113").Value = 9368.333000000001
$ws.Range("M113").Value = -4565.222
$ws.Range("N113").Value = -15876.333
$ws.Range("H116").Value = 0
$ws.Range("I116").Value = 0
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 0
$ws.Range("L116").Value = 0
$ws.Range("M116").ClearContents()
$ws.Range("H132").Value = 2180.6
$ws.Range("I132").Value = 1861.25
$ws.Range("J132").Value = 3458
$ws.Range("K132").Value = 5583.75
$ws.Range("L132").Value = 10374
$ws.Range("M132").Value = -3053.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3519.1875
$ws.Range("I32").Value = 2479.92
$ws.Range("J32").Value = 7230.857
$ws.Range("K32").Value = 2479.92
$ws.Range("L32").Value = 7230.857
$ws.Range("M32").Value = -2192.92
$ws.Range("H74").Value = 998
$ws.Range("I74").Value = 998
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 998
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -124
$ws.Range("H77").Value = 998
$ws.Range("I77").Value = 998
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 4990
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -622
$ws.Range("H97").Value = 1730.2858
$ws.Range("I97").Value = 1696.5
$ws.Range("J97").Value = 1743.8
$ws.Range("K97").Value = 1696.5
$ws.Range("L97").Value = 1743.8
$ws.Range("M97").Value = -1200.5
$ws.Range("H132").Value = 1323
$ws.Range("I132").Value = 1190.5714
$ws.Range("J132").Value = 2250
$ws.Range("K132").Value = 3571.7142
$ws.Range("L132").Value = 6750
$ws.Range("M132").Value = -1041.7142
$ws.Range("N132").Value = -11810

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H17").Value = 806.5
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 806.5
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 806.5
$ws.Range("N17").Value = -1150.5
$ws.Range("H75").Value = 15126.625
$ws.Range("I75").Value = 7004.3335
$ws.Range("J75").Value = 20000
$ws.Range("K75").Value = 7004.3335
$ws.Range("L75").Value = 20000
$ws.Range("M75").Value = -6068.3335
$ws.Range("N75").Value = -21872
$ws.Range("H78").Value = 15126.625
$ws.Range("I78").Value = 7004.3335
$ws.Range("J78").Value = 20000
$ws.Range("K78").Value = 21013.0005
$ws.Range("L78").Value = 60000
$ws.Range("M78").Value = -16333.0005
$ws.Range("N78").Value = -69360
$ws.Range("H131").Value = 60000
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 60000
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 60000
$ws.Range("N131").Value = -70080

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H29").Value = 15000
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 15000
$ws.Range("K29").Value = 0
$ws.Range("L29").Value = 15000
$ws.Range("N29").Value = -15586
$ws.Range("H31").Value = 2743.75
$ws.Range("I31").Value = 1487.5
$ws.Range("J31").Value = 4000
$ws.Range("K31").Value = 1487.5
$ws.Range("L31").Value = 4000
$ws.Range("M31").Value = -1192.5
$ws.Range("H34").Value = 2743.75
$ws.Range("I34").Value = 1487.5
$ws.Range("J34").Value = 4000
$ws.Range("K34").Value = 1487.5
$ws.Range("L34").Value = 4000
$ws.Range("M34").Value = -1285.5
$ws.Range("H57").Value = 30000000
$ws.Range("I57").Value = 0
$ws.Range("J57").Value = 30000000
$ws.Range("K57").Value = 0
$ws.Range("L57").Value = 30000000
$ws.Range("N57").Value = -30001120
$ws.Range("H107").Value = 683.3333
$ws.Range("I107").Value = 425
$ws.Range("J107").Value = 1200
$ws.Range("K107").Value = 425
$ws.Range("L107").Value = 1200
$ws.Range("M107").Value = 1495
$ws.Range("N107").Value = -5040
$ws.Range("H132").Value = 1844.0769
$ws.Range("I132").Value = 1845.88
$ws.Range("J132").Value = 1799
$ws.Range("K132").Value = 5537.64
$ws.Range("L132").Value = 5397
$ws.Range("M132").Value = -3007.64

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 290.42856
$ws.Range("I7").Value = 65.333336
$ws.Range("J7").Value = 459.25
$ws.Range("K7").Value = 196.000008
$ws.Range("L7").Value = 1377.75
$ws.Range("M7").Value = -84.00000800000001
$ws.Range("N7").Value = -1601.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 4000.1333
$ws.Range("I107").Value = 551.1818
$ws.Range("J107").Value = 13484.75
$ws.Range("K107").Value = 551.1818
$ws.Range("L107").Value = 13484.75
$ws.Range("M107").Value = 1368.8182
$ws.Range("N107").Value = -17324.75
$ws.Range("H123").Value = 39999
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 39999
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 39999
$ws.Range("N123").Value = -44899
$ws.Range("H132").Value = 1848.5
$ws.Range("I132").Value = 1848.5
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 5545.5
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -3015.5
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1117.2858
$ws.Range("I16").Value = 1136.8334
$ws.Range("J16").Value = 1000
$ws.Range("K16").Value = 1136.8334
$ws.Range("L16").Value = 1000
$ws.Range("M16").Value = -966.8334
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("M22").ClearContents()
$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("M27").ClearContents()
$ws.Range("H46").Value = 4666.6665
$ws.Range("I46").Value = 5000
$ws.Range("J46").Value = 4470.5884
$ws.Range("K46").Value = 5000
$ws.Range("L46").Value = 4470.5884
$ws.Range("M46").Value = -4812
$ws.Range("N46").Value = -4846.5884
$ws.Range("H122").Value = 1900.3334
$ws.Range("I122").Value = 1980.4
$ws.Range("J122").Value = 1500
$ws.Range("K122").Value = 5941.200000000001
$ws.Range("L122").Value = 4500
$ws.Range("M122").Value = -3491.200000000001
$ws.Range("N122").Value = -9400
$ws.Range("H132").Value = 467.8
$ws.Range("I132").Value = 467.8
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 1403.4
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = 1126.6
$ws.Range("H136").Value = 4233
$ws.Range("I136").Value = 4233
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 12699
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -10149

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 10000
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 10000
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 10000
$ws.Range("N15").Value = -10576
$ws.Range("H29").Value = 200
$ws.Range("I29").Value = 200
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 200
$ws.Range("L29").Value = 0
$ws.Range("M29").Value = 90
$ws.Range("H100").Value = 8714833
$ws.Range("I100").Value = 17425668
$ws.Range("J100").Value = 3998.75
$ws.Range("K100").Value = 34851336
$ws.Range("L100").Value = 7997.5
$ws.Range("M100").Value = -34850795
$ws.Range("N100").Value = -9079.5
$ws.Range("H126").Value = 1501.2727
$ws.Range("I126").Value = 1501.2727
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 4503.8181
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -2033.8181
$ws.Range("H132").Value = 2206.6
$ws.Range("I132").Value = 1949.9166
$ws.Range("J132").Value = 3233.3333
$ws.Range("K132").Value = 5849.7498
$ws.Range("L132").Value = 9699.999899999999
$ws.Range("M132").Value = -3319.7498
